# ITSADSSD-19935 - IAA_VevoCheck - Junk character bug fix
# Regular expression to whitelist, remove junk character &nbsp, close dashboard
# prior to logout, set time out for invokeCloseApplication workflow.
#
# Config.xlsx changes:
#  - Settings!B6 (OracleCRM_URL_DEV) no longer points at the stale
#    "uqapplications--tst1.custhelp.com" URL - it now mirrors the TEST url
#    in B7, freeing that string up for reuse.
#  - Messages sheet gains a new "RegEx_Error" message row (inserted just
#    above the existing VevoMaxCountExceed row), reusing the text that used
#    to live in Settings!B6.

$wb = $excel.ActiveWorkbook

$wsSettings = $wb.Worksheets.Item("Settings")
$wsMessages = $wb.Worksheets.Item("Messages")

# 1) Settings sheet: OracleCRM_URL_DEV (B6) now matches OracleCRM_URL_TEST (B7)
$wsSettings.Range("B6").Value = $wsSettings.Range("B7").Value2

# 2) Messages sheet: insert a new row above row 19 for the RegEx_Error message
$wsMessages.Rows("19").Insert()
$wsMessages.Range("A19").Value = "RegEx_Error"
$wsMessages.Range("B19").Value = $wsMessages.Range("B20").Value2
$wsMessages.Rows("19").RowHeight = 15

# 3) Restore selections to match the saved worksheet views
$wsMessages.Activate()
$wsMessages.Range("C27").Select()

$wsSettings.Activate()
$wsSettings.Range("A28").Select()
